# Auto-generated Excel COM-interop script
# Applies scheduled market-price refresh updates to the Ultros_Profits leve-crafting workbook
# (source: scheduled runner updating currentAveragePrice* / LevePrice* / LeveProfit* columns)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 315.07693
$ws.Range("I9").Value = 299.6
$ws.Range("K9").Value = 299.6
$ws.Range("M9").Value = -130.6
$ws.Range("H19").Value = 1562.5
$ws.Range("I19").Value = 1640
$ws.Range("J19").Value = 1433.3334
$ws.Range("K19").Value = 1640
$ws.Range("L19").Value = 1433.3334
$ws.Range("M19").Value = -1465
$ws.Range("N19").Value = -1783.3334
$ws.Range("H33").Value = 277.94736
$ws.Range("I33").Value = 292.13333
$ws.Range("K33").Value = 292.13333
$ws.Range("M33").Value = -63.13333
$ws.Range("H51").Value = 9000.666999999999
$ws.Range("I51").Value = 8000
$ws.Range("J51").Value = 9501
$ws.Range("K51").Value = 8000
$ws.Range("L51").Value = 9501
$ws.Range("M51").Value = -7516
$ws.Range("N51").Value = -10469
$ws.Range("H87").Value = 21666.445
$ws.Range("J87").Value = 21666.445
$ws.Range("L87").Value = 21666.445
$ws.Range("N87").Value = -24162.445
$ws.Range("H90").Value = 21666.445
$ws.Range("J90").Value = 21666.445
$ws.Range("L90").Value = 64999.335
$ws.Range("N90").Value = -77479.33499999999
$ws.Range("H100").Value = 7749.3
$ws.Range("I100").Value = 5898.6
$ws.Range("K100").Value = 5898.6
$ws.Range("M100").Value = -5357.6
$ws.Range("H132").Value = 13486.308
$ws.Range("I132").Value = 3448.558
$ws.Range("K132").Value = 10345.674
$ws.Range("M132").Value = -7815.673999999999
$ws.Range("H138").Value = 2670.7322
$ws.Range("I138").Value = 1189.0605
$ws.Range("J138").Value = 4796.609
$ws.Range("K138").Value = 3567.1815
$ws.Range("L138").Value = 14389.827
$ws.Range("M138").Value = 1572.8185
$ws.Range("N138").Value = -24669.827
$ws.Range("H141").Value = 2620.2173
$ws.Range("I141").Value = 2316.158
$ws.Range("J141").Value = 4064.5
$ws.Range("K141").Value = 6948.474
$ws.Range("L141").Value = 12193.5
$ws.Range("M141").Value = -1768.474
$ws.Range("N141").Value = -22553.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14497789
$ws.Range("I32").Value = 16134311
$ws.Range("K32").Value = 16134311
$ws.Range("M32").Value = -16134024
$ws.Range("H61").Value = 2866.8823
$ws.Range("I61").Value = 2312.4167
$ws.Range("K61").Value = 2312.4167
$ws.Range("M61").Value = -2100.4167
$ws.Range("H92").Value = 50930.145
$ws.Range("J92").Value = 50930.145
$ws.Range("L92").Value = 50930.145
$ws.Range("N92").Value = -55922.145
$ws.Range("H136").Value = 2866.8823
$ws.Range("I136").Value = 2312.4167
$ws.Range("K136").Value = 6937.250100000001
$ws.Range("M136").Value = -4387.250100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4262.4736
$ws.Range("I16").Value = 4381.5293
$ws.Range("J16").Value = 3250.5
$ws.Range("K16").Value = 4381.5293
$ws.Range("L16").Value = 3250.5
$ws.Range("M16").Value = -4094.5293
$ws.Range("N16").Value = -3824.5
$ws.Range("H31").Value = 2252.4102
$ws.Range("I31").Value = 2029.8518
$ws.Range("K31").Value = 2029.8518
$ws.Range("M31").Value = -1734.8518
$ws.Range("H34").Value = 2252.4102
$ws.Range("I34").Value = 2029.8518
$ws.Range("K34").Value = 2029.8518
$ws.Range("M34").Value = -1827.8518
$ws.Range("H62").Value = 10998
$ws.Range("I62").Value = 10496.5
$ws.Range("K62").Value = 10496.5
$ws.Range("M62").Value = -9872.5
$ws.Range("H65").Value = 10998
$ws.Range("I65").Value = 10496.5
$ws.Range("K65").Value = 52482.5
$ws.Range("M65").Value = -49362.5
$ws.Range("H98").Value = 31182
$ws.Range("J98").Value = 33977.5
$ws.Range("L98").Value = 33977.5
$ws.Range("N98").Value = -38469.5
$ws.Range("H99").Value = 11223471
$ws.Range("J99").Value = 18187408
$ws.Range("L99").Value = 18187408
$ws.Range("N99").Value = -18190404
$ws.Range("H113").Value = 4262.4736
$ws.Range("I113").Value = 4381.5293
$ws.Range("J113").Value = 3250.5
$ws.Range("K113").Value = 4381.5293
$ws.Range("L113").Value = 3250.5
$ws.Range("M113").Value = -2211.5293
$ws.Range("N113").Value = -7590.5
$ws.Range("H126").Value = 11223471
$ws.Range("J126").Value = 18187408
$ws.Range("L126").Value = 54562224
$ws.Range("N126").Value = -54567164
$ws.Range("H134").Value = 2927.1628
$ws.Range("I134").Value = 2721.1892
$ws.Range("J134").Value = 4197.3335
$ws.Range("K134").Value = 8163.567599999999
$ws.Range("L134").Value = 12592.0005
$ws.Range("M134").Value = -5628.567599999999
$ws.Range("N134").Value = -17662.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 111695.555
$ws.Range("I46").Value = 333766.66
$ws.Range("J46").Value = 660
$ws.Range("K46").Value = 1001299.98
$ws.Range("L46").Value = 1980
$ws.Range("M46").Value = -1001208.98
$ws.Range("N46").Value = -2162
$ws.Range("H75").Value = 1279.5714
$ws.Range("J75").Value = 2095.25
$ws.Range("L75").Value = 6285.75
$ws.Range("N75").Value = -8281.75
$ws.Range("H78").Value = 1279.5714
$ws.Range("J78").Value = 2095.25
$ws.Range("L78").Value = 18857.25
$ws.Range("N78").Value = -28841.25
$ws.Range("H92").Value = 553.5714
$ws.Range("I92").Value = 398.5
$ws.Range("J92").Value = 615.6
$ws.Range("K92").Value = 1195.5
$ws.Range("L92").Value = 1846.8
$ws.Range("M92").Value = 52.5
$ws.Range("N92").Value = -4342.8
$ws.Range("H122").Value = 737
$ws.Range("I122").Value = 749.2
$ws.Range("K122").Value = 6742.8
$ws.Range("M122").Value = -4292.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3499
$ws.Range("I126").Value = 3499
$ws.Range("K126").Value = 10497
$ws.Range("M126").Value = -8027

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4016.2856
$ws.Range("I7").Value = 2618.5
$ws.Range("K7").Value = 2618.5
$ws.Range("M7").Value = -2506.5
$ws.Range("H16").Value = 679.25
$ws.Range("I16").Value = 619.55554
$ws.Range("K16").Value = 619.55554
$ws.Range("M16").Value = -449.55554
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H40").Value = 12528
$ws.Range("I40").Value = 14861.5
$ws.Range("K40").Value = 14861.5
$ws.Range("M40").Value = -14725.5
$ws.Range("H46").Value = 8826.6
$ws.Range("I46").Value = 3700
$ws.Range("K46").Value = 3700
$ws.Range("M46").Value = -3512
$ws.Range("H51").Value = 30247.5
$ws.Range("J51").Value = 30247.5
$ws.Range("L51").Value = 30247.5
$ws.Range("N51").Value = -31203.5
$ws.Range("H55").Value = 1660
$ws.Range("I55").Value = 1403.4445
$ws.Range("J55").Value = 1916.5555
$ws.Range("K55").Value = 1403.4445
$ws.Range("L55").Value = 1916.5555
$ws.Range("M55").Value = -1230.4445
$ws.Range("N55").Value = -2262.5555
$ws.Range("H126").Value = 4016.2856
$ws.Range("I126").Value = 2618.5
$ws.Range("K126").Value = 7855.5
$ws.Range("M126").Value = -5385.5
$ws.Range("H136").Value = 2022.8689
$ws.Range("I136").Value = 1816.017
$ws.Range("K136").Value = 5448.051
$ws.Range("M136").Value = -2898.051

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 12250.2
$ws.Range("I74").Value = 10000
$ws.Range("J74").Value = 12812.75
$ws.Range("K74").Value = 10000
$ws.Range("L74").Value = 12812.75
$ws.Range("M74").Value = -9064
$ws.Range("N74").Value = -14684.75
$ws.Range("H77").Value = 12250.2
$ws.Range("I77").Value = 10000
$ws.Range("J77").Value = 12812.75
$ws.Range("K77").Value = 30000
$ws.Range("L77").Value = 38438.25
$ws.Range("M77").Value = -25320
$ws.Range("N77").Value = -47798.25
$ws.Range("H96").Value = 107043.4
$ws.Range("I96").Value = 129304.25
$ws.Range("K96").Value = 129304.25
$ws.Range("M96").Value = -127931.25
$ws.Range("H126").Value = 2119.158
$ws.Range("I126").Value = 1325.7
$ws.Range("K126").Value = 3977.1
$ws.Range("M126").Value = -1507.1
$ws.Range("H132").Value = 1651.9706
$ws.Range("I132").Value = 1380.25
$ws.Range("K132").Value = 4140.75
$ws.Range("M132").Value = -1610.75
$ws.Range("H136").Value = 1626.6666
$ws.Range("I136").Value = 956.8
$ws.Range("J136").Value = 2743.111
$ws.Range("K136").Value = 2870.4
$ws.Range("L136").Value = 8229.332999999999
$ws.Range("M136").Value = -320.3999999999996
$ws.Range("N136").Value = -13329.333
